$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "base de datos" (Periodo Mora table, rows 16-18) is reshuffled:
# a new "2402" period row for JAVIER EDUARDO PUERTA COLINA is inserted first,
# the existing "2401" row for JAVIER EDUARDO PUERTA COLINA moves to row 17,
# and the "2401" row for DOVANIS DE JESUS MONTAÑO VIADERO moves to row 18.

# Row 16: PPT / 7448501 / JAVIER EDUARDO PUERTA COLINA / 2402 / 1600 / 1200000
$ws.Range("B16").Value = "PPT"
$ws.Range("C16").Value = "7448501"
$ws.Range("D16").Value = "JAVIER EDUARDO PUERTA COLINA"
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 1600
$ws.Range("G16").Value = 1200000

# Row 17: PPT / 7448501 / JAVIER EDUARDO PUERTA COLINA / 2401 / 48000 / 1200000
$ws.Range("B17").Value = "PPT"
$ws.Range("C17").Value = "7448501"
$ws.Range("D17").Value = "JAVIER EDUARDO PUERTA COLINA"
$ws.Range("E17").Value = "2401"
$ws.Range("F17").Value = 48000
$ws.Range("G17").Value = 1200000

# Row 18: CC / 1127591543 / DOVANIS DE JESUS MONTAÑO VIADERO / 2401 / 6933 / 1300000
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1127591543"
$ws.Range("D18").Value = "DOVANIS DE JESUS MONTAÑO VIADERO"
$ws.Range("E18").Value = "2401"
$ws.Range("F18").Value = 6933
$ws.Range("G18").Value = 1300000
